$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Fill in the week-7 (column U) and week-8 (column X) hours that were
# previously left blank for Research (row 6), Thinking (row 8) and
# Programming (row 12).
$ws.Range("U6").Value = 5
$ws.Range("X6").Value = 2

$ws.Range("U8").Value = 4
$ws.Range("X8").Value = 2

$ws.Range("U12").Value = 5
$ws.Range("X12").Value = 1

# The week-7 total (U20) now derives from the column instead of being a
# hard-coded 0, matching the live SUM formulas already used by the other
# week totals (O20, R20).
$ws.Range("U20").Formula = "=SUM(U6:W19)"

# Scroll the sheet down one row and move the selection onto the
# Programming/week-6 block.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("R12:T12").Select()
